# Apply "Trade #50 closed at 2026-02-17 21:08:10" update to live_trading_results.xlsx
$wb = $excel.ActiveWorkbook

$wsSummary   = $wb.Worksheets.Item("Summary")
$wsStrategy  = $wb.Worksheets.Item("Strategy Status")
$wsAllTrades = $wb.Worksheets.Item("All Trades")
$wsMM        = $wb.Worksheets.Item("MarketMaking")

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$wsSummary.Range("B3").Value = 1400.6    # Current Capital
$wsSummary.Range("B4").Value = 0.4       # Total P&L $
$wsSummary.Range("B6").Value = 78        # Total Trades
$wsSummary.Range("B7").Value = 37        # Winning Trades
$wsSummary.Range("B9").Value = 47.44     # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 5)
# ---------------------------------------------------------------------------
$wsStrategy.Range("C5").Value = 100.6    # Capital
$wsStrategy.Range("D5").Value = 45       # Trades
$wsStrategy.Range("E5").Value = 0.29     # P&L $
$wsStrategy.Range("F5").Value = 0.6      # P&L %
$wsStrategy.Range("G5").Value = 51.11    # Win Rate %

# ---------------------------------------------------------------------------
# All Trades sheet - update existing Trade #78 row (row 79) now that it closed
# ---------------------------------------------------------------------------
$wsAllTrades.Cells.Item(79, 7).Value = 0.97               # G: Exit Price
$wsAllTrades.Cells.Item(79, 8).Value = "CLOSED"            # H: Status
$wsAllTrades.Cells.Item(79, 9).Value = 2.1053              # I: P&L %
$wsAllTrades.Cells.Item(79, 10).Value = 0.02                # J: P&L $
$wsAllTrades.Cells.Item(79, 11).Value = 100.6               # K: Capital After
$wsAllTrades.Cells.Item(79, 12).Value = "early_exit"        # L: Exit Reason
$wsAllTrades.Cells.Item(79, 13).Value = 0.14                # M: Duration (min)

# All Trades sheet - append new Trade #111 row (row 112)
$wsAllTrades.Cells.Item(112, 1).Value = 111
$wsAllTrades.Cells.Item(112, 2).Value = "'2026-02-17"
$wsAllTrades.Cells.Item(112, 3).Value = "21:08:04"
$wsAllTrades.Cells.Item(112, 4).Value = "MarketMaking"
$wsAllTrades.Cells.Item(112, 5).Value = "UP"
$wsAllTrades.Cells.Item(112, 6).Value = 0.95
$wsAllTrades.Cells.Item(112, 8).Value = "OPEN"
$wsAllTrades.Cells.Item(112, 9).Value = 0
$wsAllTrades.Cells.Item(112, 10).Value = 0
$wsAllTrades.Cells.Item(112, 11).Value = 100.5814872031006
$wsAllTrades.Cells.Item(112, 13).Value = 0
$wsAllTrades.Cells.Item(112, 14).Value = 0
$wsAllTrades.Cells.Item(112, 15).Value = 0
$wsAllTrades.Cells.Item(112, 16).Value = 0.6
$wsAllTrades.Cells.Item(112, 17).Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------------
# MarketMaking sheet - update existing Trade #78 row (row 46) now that it closed
# ---------------------------------------------------------------------------
$wsMM.Cells.Item(46, 7).Value = 0.97                # G: Exit Price
$wsMM.Cells.Item(46, 8).Value = "CLOSED"             # H: Status
$wsMM.Cells.Item(46, 9).Value = 2.1053               # I: P&L %
$wsMM.Cells.Item(46, 10).Value = 0.02                 # J: P&L $
$wsMM.Cells.Item(46, 11).Value = 100.6                # K: Capital After
$wsMM.Cells.Item(46, 16).Value = "early_exit"         # P: Exit Reason
$wsMM.Cells.Item(46, 17).Value = 0.14                 # Q: Duration (min)

# MarketMaking sheet - append new Trade #111 row (row 79)
$wsMM.Cells.Item(79, 1).Value = 111
$wsMM.Cells.Item(79, 2).Value = "'2026-02-17"
$wsMM.Cells.Item(79, 3).Value = "21:08:04"
$wsMM.Cells.Item(79, 4).Value = "MarketMaking"
$wsMM.Cells.Item(79, 5).Value = "UP"
$wsMM.Cells.Item(79, 6).Value = 0.95
$wsMM.Cells.Item(79, 8).Value = "OPEN"
$wsMM.Cells.Item(79, 9).Value = 0
$wsMM.Cells.Item(79, 10).Value = 0
$wsMM.Cells.Item(79, 11).Value = 100.5814872031006
$wsMM.Cells.Item(79, 12).Value = 0
$wsMM.Cells.Item(79, 13).Value = 0
$wsMM.Cells.Item(79, 14).Value = 0.6
$wsMM.Cells.Item(79, 15).Value = "Normal spread capture: 19600 bps"
$wsMM.Cells.Item(79, 17).Value = 0

Write-Host "Edit complete"
